$wb = $excel.ActiveWorkbook

# --- Sheet: 1er Parcial ---
$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("C8").Value = 81
$ws.Range("D8").Value = 48.5
$ws.Range("E8").Value = 86
$ws.Range("F8").Value = 51.5
$ws.Range("I8").Value = 7.9

$ws.Range("C9").Value = 46
$ws.Range("D9").Value = 47.42
$ws.Range("E9").Value = 51
$ws.Range("F9").Value = 52.58
$ws.Range("I9").Value = 9.1

$ws.Range("C10").Value = 51
$ws.Range("D10").Value = 29.65
$ws.Range("E10").Value = 121
$ws.Range("F10").Value = 70.34999999999999
$ws.Range("I10").Value = 8.300000000000001

$ws.Range("C11").Value = 29
$ws.Range("D11").Value = 28.43
$ws.Range("E11").Value = 73
$ws.Range("F11").Value = 71.56999999999999
$ws.Range("I11").Value = 8.300000000000001

$ws.Range("C13").Value = 42
$ws.Range("D13").Value = 26.42
$ws.Range("E13").Value = 111
$ws.Range("F13").Value = 69.81
$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 3.77
$ws.Range("I13").Value = 8

$ws.Range("C14").Value = 23
$ws.Range("D14").Value = 12.17
$ws.Range("E14").Value = 147
$ws.Range("F14").Value = 77.78
$ws.Range("G14").Value = 19
$ws.Range("H14").Value = 10.05
$ws.Range("I14").Value = 8

$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 69
$ws.Range("F15").Value = 71.88
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = 28.13
$ws.Range("I15").Value = 7.2

$ws.Range("C16").Value = 35
$ws.Range("D16").Value = 19.66
$ws.Range("E16").Value = 143
$ws.Range("F16").Value = 80.34
$ws.Range("I16").Value = 8.300000000000001

$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 19.73
$ws.Range("E18").Value = 118
$ws.Range("F18").Value = 80.27

$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 9.210000000000001
$ws.Range("E19").Value = 164
$ws.Range("F19").Value = 71.93000000000001
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 18.86
$ws.Range("I19").Value = 7.6

$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 35.83
$ws.Range("E21").Value = 77
$ws.Range("F21").Value = 64.17
$ws.Range("I21").Value = 8.5

$ws.Range("C22").Value = 48
$ws.Range("D22").Value = 22.64
$ws.Range("E22").Value = 164
$ws.Range("F22").Value = 77.36
$ws.Range("I22").Value = 7.4

$ws.Range("C23").Value = 56
$ws.Range("D23").Value = 28.72
$ws.Range("E23").Value = 139
$ws.Range("F23").Value = 71.28
$ws.Range("I23").Value = 9.199999999999999

$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 67
$ws.Range("F25").Value = 64.42
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 35.58
$ws.Range("I25").Value = 6.4

$ws.Range("C26").Value = 61
$ws.Range("D26").Value = 58.1
$ws.Range("E26").Value = 44
$ws.Range("F26").Value = 41.9
$ws.Range("I26").Value = 8.6

$ws.Range("C28").Value = 35
$ws.Range("D28").Value = 25.74
$ws.Range("E28").Value = 101
$ws.Range("F28").Value = 74.26000000000001
$ws.Range("I28").Value = 10

$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 24
$ws.Range("F30").Value = 100
$ws.Range("I30").Value = 8.800000000000001

$ws.Range("C31").Value = 32
$ws.Range("D31").Value = 29.63
$ws.Range("E31").Value = 76
$ws.Range("F31").Value = 70.37
$ws.Range("I31").Value = 8.199999999999999

$ws.Range("C32").Value = 102
$ws.Range("D32").Value = 41.3
$ws.Range("E32").Value = 145
$ws.Range("F32").Value = 58.7
$ws.Range("I32").Value = 7.9

$ws.Range("C35").Value = 143
$ws.Range("D35").Value = 88.27
$ws.Range("E35").Value = 19
$ws.Range("F35").Value = 11.73
$ws.Range("I35").Value = 9.300000000000001

$ws.Range("C36").Value = 39
$ws.Range("D36").Value = 27.46
$ws.Range("E36").Value = 98
$ws.Range("F36").Value = 69.01000000000001
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 3.52
$ws.Range("I36").Value = 7.6

$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 65
$ws.Range("F37").Value = 100
$ws.Range("I37").Value = 8.1

$ws.Range("C39").Value = 96
$ws.Range("D39").Value = 68.09
$ws.Range("E39").Value = 45
$ws.Range("F39").Value = 31.91
$ws.Range("I39").Value = 7.9

$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 0.98
$ws.Range("E40").Value = 167
$ws.Range("F40").Value = 81.86
$ws.Range("G40").Value = 35
$ws.Range("H40").Value = 17.16
$ws.Range("I40").Value = 7.1

$ws.Range("C41").Value = 23
$ws.Range("D41").Value = 31.94
$ws.Range("E41").Value = 49
$ws.Range("F41").Value = 68.06
$ws.Range("I41").Value = 7.5

$ws.Range("C43").Value = 14
$ws.Range("D43").Value = 42.42
$ws.Range("E43").Value = 19
$ws.Range("F43").Value = 57.58
$ws.Range("I43").Value = 7.4

$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 108
$ws.Range("F47").Value = 69.23
$ws.Range("G47").Value = 48
$ws.Range("H47").Value = 30.77

$ws.Range("C48").Value = 44
$ws.Range("D48").Value = 30.56
$ws.Range("E48").Value = 100
$ws.Range("F48").Value = 69.44
$ws.Range("I48").Value = 6.8

$ws.Range("C49").Value = 59
$ws.Range("D49").Value = 47.58
$ws.Range("E49").Value = 65
$ws.Range("F49").Value = 52.42
$ws.Range("I49").Value = 8.800000000000001


# --- Sheet: 2o Parcial ---
$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("G8").Value = 86
$ws.Range("H8").Value = 51.5

$ws.Range("G9").Value = 51
$ws.Range("H9").Value = 52.58

$ws.Range("G10").Value = 121
$ws.Range("H10").Value = 70.34999999999999

$ws.Range("G11").Value = 73
$ws.Range("H11").Value = 71.56999999999999

$ws.Range("G13").Value = 117
$ws.Range("H13").Value = 73.58

$ws.Range("G14").Value = 166
$ws.Range("H14").Value = 87.83

$ws.Range("G15").Value = 96
$ws.Range("H15").Value = 100

$ws.Range("G16").Value = 143
$ws.Range("H16").Value = 80.34

$ws.Range("G18").Value = 118
$ws.Range("H18").Value = 80.27

$ws.Range("G19").Value = 207
$ws.Range("H19").Value = 90.79000000000001

$ws.Range("G21").Value = 77
$ws.Range("H21").Value = 64.17

$ws.Range("G22").Value = 164
$ws.Range("H22").Value = 77.36

$ws.Range("G23").Value = 139
$ws.Range("H23").Value = 71.28

$ws.Range("G25").Value = 104
$ws.Range("H25").Value = 100

$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 41.9

$ws.Range("G28").Value = 101
$ws.Range("H28").Value = 74.26000000000001

$ws.Range("G30").Value = 24
$ws.Range("H30").Value = 100

$ws.Range("G31").Value = 76
$ws.Range("H31").Value = 70.37

$ws.Range("G32").Value = 145
$ws.Range("H32").Value = 58.7

$ws.Range("G35").Value = 19
$ws.Range("H35").Value = 11.73

$ws.Range("G36").Value = 103
$ws.Range("H36").Value = 72.54000000000001

$ws.Range("G37").Value = 65
$ws.Range("H37").Value = 100

$ws.Range("G39").Value = 45
$ws.Range("H39").Value = 31.91

$ws.Range("G40").Value = 202
$ws.Range("H40").Value = 99.02

$ws.Range("G41").Value = 49
$ws.Range("H41").Value = 68.06

$ws.Range("G43").Value = 19
$ws.Range("H43").Value = 57.58

$ws.Range("G47").Value = 156
$ws.Range("H47").Value = 100

$ws.Range("G48").Value = 100
$ws.Range("H48").Value = 69.44

$ws.Range("G49").Value = 65
$ws.Range("H49").Value = 52.42


# --- Sheet: Final ---
$ws = $wb.Worksheets.Item("Final")
$ws.Range("C8").Value = 81
$ws.Range("D8").Value = 48.5
$ws.Range("E8").Value = 86
$ws.Range("F8").Value = 51.5
$ws.Range("I8").Value = 7.9

$ws.Range("C9").Value = 46
$ws.Range("D9").Value = 47.42
$ws.Range("E9").Value = 51
$ws.Range("F9").Value = 52.58
$ws.Range("I9").Value = 9.1

$ws.Range("C10").Value = 51
$ws.Range("D10").Value = 29.65
$ws.Range("E10").Value = 121
$ws.Range("F10").Value = 70.34999999999999
$ws.Range("I10").Value = 8.300000000000001

$ws.Range("C11").Value = 29
$ws.Range("D11").Value = 28.43
$ws.Range("E11").Value = 73
$ws.Range("F11").Value = 71.56999999999999
$ws.Range("I11").Value = 8.300000000000001

$ws.Range("C13").Value = 42
$ws.Range("D13").Value = 26.42
$ws.Range("E13").Value = 111
$ws.Range("F13").Value = 69.81
$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 3.77
$ws.Range("I13").Value = 8

$ws.Range("C14").Value = 23
$ws.Range("D14").Value = 12.17
$ws.Range("E14").Value = 147
$ws.Range("F14").Value = 77.78
$ws.Range("G14").Value = 19
$ws.Range("H14").Value = 10.05
$ws.Range("I14").Value = 8

$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 69
$ws.Range("F15").Value = 71.88
$ws.Range("G15").Value = 27
$ws.Range("H15").Value = 28.13
$ws.Range("I15").Value = 7.2

$ws.Range("C16").Value = 35
$ws.Range("D16").Value = 19.66
$ws.Range("E16").Value = 143
$ws.Range("F16").Value = 80.34
$ws.Range("I16").Value = 8.300000000000001

$ws.Range("C18").Value = 29
$ws.Range("D18").Value = 19.73
$ws.Range("E18").Value = 118
$ws.Range("F18").Value = 80.27

$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 9.210000000000001
$ws.Range("E19").Value = 164
$ws.Range("F19").Value = 71.93000000000001
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 18.86
$ws.Range("I19").Value = 7.6

$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 35.83
$ws.Range("E21").Value = 77
$ws.Range("F21").Value = 64.17
$ws.Range("I21").Value = 8.5

$ws.Range("C22").Value = 48
$ws.Range("D22").Value = 22.64
$ws.Range("E22").Value = 164
$ws.Range("F22").Value = 77.36
$ws.Range("I22").Value = 7.4

$ws.Range("C23").Value = 56
$ws.Range("D23").Value = 28.72
$ws.Range("E23").Value = 139
$ws.Range("F23").Value = 71.28
$ws.Range("I23").Value = 9.199999999999999

$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 67
$ws.Range("F25").Value = 64.42
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 35.58
$ws.Range("I25").Value = 6.4

$ws.Range("C26").Value = 61
$ws.Range("D26").Value = 58.1
$ws.Range("E26").Value = 44
$ws.Range("F26").Value = 41.9
$ws.Range("I26").Value = 8.6

$ws.Range("C28").Value = 35
$ws.Range("D28").Value = 25.74
$ws.Range("E28").Value = 101
$ws.Range("F28").Value = 74.26000000000001
$ws.Range("I28").Value = 10

$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 24
$ws.Range("F30").Value = 100
$ws.Range("I30").Value = 8.800000000000001

$ws.Range("C31").Value = 32
$ws.Range("D31").Value = 29.63
$ws.Range("E31").Value = 76
$ws.Range("F31").Value = 70.37
$ws.Range("I31").Value = 8.199999999999999

$ws.Range("C32").Value = 102
$ws.Range("D32").Value = 41.3
$ws.Range("E32").Value = 145
$ws.Range("F32").Value = 58.7
$ws.Range("I32").Value = 7.9

$ws.Range("C35").Value = 143
$ws.Range("D35").Value = 88.27
$ws.Range("E35").Value = 19
$ws.Range("F35").Value = 11.73
$ws.Range("I35").Value = 9.300000000000001

$ws.Range("C36").Value = 39
$ws.Range("D36").Value = 27.46
$ws.Range("E36").Value = 98
$ws.Range("F36").Value = 69.01000000000001
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 3.52
$ws.Range("I36").Value = 7.6

$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 65
$ws.Range("F37").Value = 100
$ws.Range("I37").Value = 8.1

$ws.Range("C39").Value = 96
$ws.Range("D39").Value = 68.09
$ws.Range("E39").Value = 45
$ws.Range("F39").Value = 31.91
$ws.Range("I39").Value = 7.9

$ws.Range("C40").Value = 2
$ws.Range("D40").Value = 0.98
$ws.Range("E40").Value = 167
$ws.Range("F40").Value = 81.86
$ws.Range("G40").Value = 35
$ws.Range("H40").Value = 17.16
$ws.Range("I40").Value = 7.1

$ws.Range("C41").Value = 23
$ws.Range("D41").Value = 31.94
$ws.Range("E41").Value = 49
$ws.Range("F41").Value = 68.06
$ws.Range("I41").Value = 7.5

$ws.Range("C43").Value = 14
$ws.Range("D43").Value = 42.42
$ws.Range("E43").Value = 19
$ws.Range("F43").Value = 57.58
$ws.Range("I43").Value = 7.4

$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 108
$ws.Range("F47").Value = 69.23
$ws.Range("G47").Value = 48
$ws.Range("H47").Value = 30.77

$ws.Range("C48").Value = 44
$ws.Range("D48").Value = 30.56
$ws.Range("E48").Value = 100
$ws.Range("F48").Value = 69.44
$ws.Range("I48").Value = 6.8

$ws.Range("C49").Value = 59
$ws.Range("D49").Value = 47.58
$ws.Range("E49").Value = 65
$ws.Range("F49").Value = 52.42
$ws.Range("I49").Value = 8.800000000000001

